$wb = $excel.ActiveWorkbook

# --- Sheet "About" (sheet1) ---
$about = $wb.Worksheets.Item("About")

# Remove the "-> electricity source within all fuels" row: it used to sit at A10
# (between "hydrogen combined cycle es" at row 9-ish context and "hard coal es").
# After the edit, the whole list A10:A33 is shifted up by one row, and the
# previously-last row (A34, "hydrogen combined cycle es") is dropped since the
# list now ends at A33 with "hydrogen combined cycle es" (previously at A33 "-> ..").
# Simplest: delete the row that currently contains that text, which shifts
# everything up automatically.
$found = $about.Columns.Item(1).Find("-> electricity source within all fuels")
if ($found -ne $null) {
    $found.EntireRow.Delete() | Out-Null
}

# --- Sheet "ESUfR" (sheet2) ---
$esufr = $wb.Worksheets.Item("ESUfR")

# Capture the old "X es" strings from column A (rows 2-9) before we overwrite them.
$oldValues = @()
for ($r = 2; $r -le 9; $r++) {
    $oldValues += $esufr.Cells.Item($r, 1).Value2
}

# Move the old header "Electricity Sources" from A1 to B1, and set new A1 header.
$oldHeader = $esufr.Cells.Item(1, 1).Value2
$esufr.Cells.Item(1, 2).Value = $oldHeader
$esufr.Cells.Item(1, 1).Value = "Electricity Sources (no es)"
$esufr.Cells.Item(1, 2).Font.Bold = $true

# Re-write column A with the "no es" names (strip trailing " es").
for ($i = 0; $i -lt $oldValues.Length; $i++) {
    $r = $i + 2
    $base = $oldValues[$i] -replace ' es$', ''
    $esufr.Cells.Item($r, 1).Value = $base
}

# Column B rows 2-21: formula that re-appends " es" (shared-formula equivalent).
# B2 is written individually and B3:B21 as one range so the resulting shared
# formula group spans B3:B21 (matching the authored layout).
$esufr.Range("B2").Formula = '=IF(A2="","",CONCATENATE(A2," es"))'
$esufr.Range("B3:B21").Formula = '=IF(A3="","",CONCATENATE(A3," es"))'

$esufr.Range("C45").Select() | Out-Null

# Leave the "About" sheet selected/active at the end, matching the
# final author selection state.
$about.Range("D53").Select() | Out-Null
